$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap row 28 and row 29 data (per-column, only columns that differ)
$ws.Range("A28").Value = 111596831
$ws.Range("A29").Value = 111596897
$ws.Range("B28").Value = 89405
$ws.Range("B29").Value = 103288
$ws.Range("D28").Value = 'NT'
$ws.Range("D29").Value = 'LC'
$ws.Range("E28").Value = 1202
$ws.Range("E29").Value = 221144
$ws.Range("F28").Value = 'Ullticka'
$ws.Range("F29").Value = 'Grönpyrola'
$ws.Range("G28").Value = 'Phellinidium ferrugineofuscum'
$ws.Range("G29").Value = 'Pyrola chlorantha'
$ws.Range("H28").Value = '(P.Karst.) Fiasson & Niemelä'
$ws.Range("H29").Value = 'Sw.'
$ws.Range("I28").Value = $null
$ws.Range("I29").Value = '''100'
$ws.Range("J28").Value = $null
$ws.Range("J29").Value = 'plantor/tuvor'
$ws.Range("P28").Value = 'Björkmossen 304 m N, Upl'
$ws.Range("P29").Value = 'Björkmossen 227 m E, Upl'
$ws.Range("Q28").Value = 654194.6095515667
$ws.Range("Q29").Value = 654422.181084068
$ws.Range("R28").Value = 6691076.210478476
$ws.Range("R29").Value = 6690769.97221576
$ws.Range("S28").Value = 4
$ws.Range("S29").Value = 8
$ws.Range("Z28").Value = '16:20'
$ws.Range("Z29").Value = '12:53'
$ws.Range("AB28").Value = '16:20'
$ws.Range("AB29").Value = '12:53'
$ws.Range("AC28").Value = $null
$ws.Range("AC29").Value = 'Uppskattat antal, helt tjockt med plantor så går ej att räkna.'

# Swap row 38 and row 39 data (per-column, only columns that differ)
$ws.Range("A38").Value = 111596895
$ws.Range("A39").Value = 111596859
$ws.Range("B38").Value = 103288
$ws.Range("B39").Value = 99413
$ws.Range("E38").Value = 221144
$ws.Range("E39").Value = 221235
$ws.Range("F38").Value = 'Grönpyrola'
$ws.Range("F39").Value = 'Vårärt'
$ws.Range("G38").Value = 'Pyrola chlorantha'
$ws.Range("G39").Value = 'Lathyrus vernus'
$ws.Range("H38").Value = 'Sw.'
$ws.Range("H39").Value = '(L.) Bernh.'
$ws.Range("I38").Value = '''30'
$ws.Range("I39").Value = '''15'
$ws.Range("P38").Value = 'Björkmossen 238 m E, Upl'
$ws.Range("P39").Value = 'Björkmossen 311 m SE, Upl'
$ws.Range("Q38").Value = 654433.1528313066
$ws.Range("Q39").Value = 654503.0812791266
$ws.Range("R38").Value = 6690768.95009726
$ws.Range("R39").Value = 6690724.805518131
$ws.Range("Z38").Value = '12:48'
$ws.Range("Z39").Value = '11:25'
$ws.Range("AB38").Value = '12:48'
$ws.Range("AB39").Value = '11:27'
$ws.Range("AC38").Value = 'Cirka.'
$ws.Range("AC39").Value = $null
